$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D19").Value = "    "
$ws.Range("E19").Select()
